$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.964999999999999
$ws.Range("A3").Value = -21.663
$ws.Range("C3").Value = -12.622
$ws.Range("C12").Value = -11.536
$ws.Range("A14").Value = -21.682
$ws.Range("A16").Value = -21.373
$ws.Range("B18").Value = 5.626
$ws.Range("A21").Value = -20.421
$ws.Range("A23").Value = -20.428
$ws.Range("B24").Value = 6.879
$ws.Range("C24").Value = -13.45
$ws.Range("A25").Value = -20.837
$ws.Range("B25").Value = 6.815
$ws.Range("C25").Value = -12.961
$ws.Range("A26").Value = -21.374
$ws.Range("B27").Value = 5.513
$ws.Range("A29").Value = -21.202
$ws.Range("B30").Value = 6.245
$ws.Range("B31").Value = 6.117000000000001
$ws.Range("B39").Value = 7.722
$ws.Range("A40").Value = -20.294
$ws.Range("C41").Value = -12.61
$ws.Range("B42").Value = 8.494
$ws.Range("B48").Value = 5.24
$ws.Range("C50").Value = -13.087
$ws.Range("B51").Value = 5.330999999999999
$ws.Range("B52").Value = 5.415999999999999
$ws.Range("A53").Value = -21.793
$ws.Range("C53").Value = -11.726
$ws.Range("B55").Value = 4.703999999999999
$ws.Range("B56").Value = 5.801
$ws.Range("C56").Value = -12.989
$ws.Range("A57").Value = -21.352
$ws.Range("B57").Value = 6.095000000000001
$ws.Range("C57").Value = -13.342
$ws.Range("C58").Value = -13.081
$ws.Range("A59").Value = -22.404
$ws.Range("B60").Value = 5.896999999999999
$ws.Range("C61").Value = -13.181
$ws.Range("C63").Value = -11.864
$ws.Range("C64").Value = -11.684
$ws.Range("A65").Value = -21.329
$ws.Range("A69").Value = -21.507
$ws.Range("C70").Value = -11.536
$ws.Range("C72").Value = -11.822
$ws.Range("B73").Value = 6.187
$ws.Range("B74").Value = 8.403
$ws.Range("A79").Value = -21.251
$ws.Range("A83").Value = -21.212
$ws.Range("C86").Value = -13.036
$ws.Range("B89").Value = 5.549000000000001
$ws.Range("C89").Value = -11.134
$ws.Range("B90").Value = 5.906999999999999
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 6.027
$ws.Range("A93").Value = -21.324
$ws.Range("C98").Value = -12.39
$ws.Range("A100").Value = -21.424
$ws.Range("C100").Value = -13.108
$ws.Range("C102").Value = -13.237
